$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

# The target values collide with other cells' original/intermediate values
# (e.g. row 8's new AIC "2,138.00" equals row 10's original AIC), so a
# direct single-pass old->new replacement can retarget the wrong cell.
# Work around it with a two-phase placeholder swap: every "old" string
# below is unique in the starting document, so phase 1 can't collide;
# every placeholder is unique and unambiguous, so phase 2 can't either.

# ---- Phase 1: old text -> unique placeholder ----
Replace-Exact "RT ~ condition_block + trial_block + CL + SEval + (CL + SEval|ID)" "@@P1@@"
Replace-Exact "1,862.00" "@@P2@@"
Replace-Exact "1,926.00" "@@P3@@"

Replace-Exact "RT ~ condition_block + trial_block + CL * SEval + (CL + SEval|ID)" "@@P4@@"
Replace-Exact "1,865.00" "@@P5@@"
Replace-Exact "1,934.00" "@@P6@@"

Replace-Exact "RT ~ condition_block + CL + SEval + (CL + SEval|ID)" "@@P7@@"
Replace-Exact "1,943.00" "@@P8@@"
Replace-Exact "1,998.00" "@@P9@@"

Replace-Exact "RT ~ trial_block + CL + SEval + (CL + SEval|ID)" "@@P10@@"
Replace-Exact "1,982.00" "@@P11@@"
Replace-Exact "2,041.00" "@@P12@@"

Replace-Exact "RT ~ CL + SEval + (CL + SEval|ID)" "@@P13@@"
Replace-Exact "2,051.00" "@@P14@@"
Replace-Exact "2,100.00" "@@P15@@"

Replace-Exact "RT ~ CL * SEval + (CL + SEval|ID)" "@@P16@@"
Replace-Exact "2,055.00" "@@P17@@"
Replace-Exact "2,109.00" "@@P18@@"

Replace-Exact "RT ~ CL + SEval + (1|ID)" "@@P19@@"
Replace-Exact "2,138.00" "@@P20@@"
Replace-Exact "2,163.00" "@@P21@@"

# ---- Phase 2: placeholder -> final text ----
Replace-Exact "@@P1@@" "RT ~ condition_block + trial_block + CL + SEval + (1|ID)"
Replace-Exact "@@P2@@" "1,990.00"
Replace-Exact "@@P3@@" "2,029.00"

Replace-Exact "@@P4@@" "RT ~ condition_block + trial_block + CL * SEval + (1|ID)"
Replace-Exact "@@P5@@" "1,993.00"
Replace-Exact "@@P6@@" "2,038.00"

Replace-Exact "@@P7@@" "RT ~ condition_block + CL + SEval + (1|ID)"
Replace-Exact "@@P8@@" "2,051.00"
Replace-Exact "@@P9@@" "2,081.00"

Replace-Exact "@@P10@@" "RT ~ trial_block + CL + SEval + (1|ID)"
Replace-Exact "@@P11@@" "2,083.00"
Replace-Exact "@@P12@@" "2,118.00"

Replace-Exact "@@P13@@" "RT ~ CL + SEval + (1|ID)"
Replace-Exact "@@P14@@" "2,138.00"
Replace-Exact "@@P15@@" "2,163.00"

Replace-Exact "@@P16@@" "RT ~ CL + SEval + (1|ID)"
Replace-Exact "@@P17@@" "2,142.00"
Replace-Exact "@@P18@@" "2,171.00"

Replace-Exact "@@P19@@" "RT ~ CL * SEval + (1|ID)"
Replace-Exact "@@P20@@" "2,142.00"
Replace-Exact "@@P21@@" "2,171.00"
